# Edit summary (per the target diff):
#   1. On slide 16, the table's style (tableStyleId) changes from
#      {8ADB841B-590F-412E-B102-28B5B76EE8B4} to
#      {4B286482-5C10-480D-BC3C-6EDA48513499}.
#   2. The deck's theme (ppt/theme/theme1.xml, bound to the one slide
#      master used by every layout/slide) is switched from the
#      "Integral" palette to the standard Office palette.
#
# Note: the authored diff also shows ppt/theme/theme2.xml (consumed only
# by the Notes Master, not by any visible slide) swapping the other way
# (Office -> Integral). PowerPoint's object model only ever exposes a
# single Theme/ThemeColorScheme for the deck (Master.Theme,
# NotesMaster.Theme, and HandoutMaster.Theme all resolve to the same
# object backed by the slide master's theme part), so the Notes-Master
# theme part is not independently reachable here, and the clrScheme/
# theme "name" attributes are read-only. Everything that is reachable
# through the PowerPoint COM surface is applied below.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{4B286482-5C10-480D-BC3C-6EDA48513499}")

# --- 2. Theme colors: Integral -> Office -----------------------------
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72
